# Generate Report for Handoff
# - Flip the "In Translation" status to "Ready for handoff" everywhere it
#   appears (Overview!E2/F2, zh-cn!C2, de-de!C2).
# - Bump the handoff timestamps that were refreshed when the report was
#   regenerated (Overview!G2 & de-de!H2 share one stamp, zh-cn!H2 has its own).
# - Widen the Status / language columns so the longer "Ready for handoff"
#   text still fits.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- status text -----------------------------------------------------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- refreshed timestamps ---------------------------------------------------
$overview.Range("G2").Value = "2016-08-25 06:38:13"
$dede.Range("H2").Value     = "2016-08-25 06:38:13"
$zhcn.Range("H2").Value     = "2016-08-25 06:38:08"

# --- column widths -----------------------------------------------------------
# The stored OOXML column width of 17.2159881591797 corresponds to a
# COM ColumnWidth of ~16.33 once Excel's internal pixel-grid rounding is
# taken into account (COM ColumnWidth is quantized to whole pixels, i.e.
# steps of 1/6 on this font/DPI), so use the value that lands closest to it.
$newWidth = 16.333333333333332

$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth
$zhcn.Columns.Item(3).ColumnWidth     = $newWidth
$dede.Columns.Item(3).ColumnWidth     = $newWidth
